$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# FC_WEGE3: insert a new "CAPEX" column before column V (old V..AN
# shift right to W..AO). New column V gets a formula that sums the
# (now shifted) "Compra de invest perman" (X) and "Compra ativos fix
# e dif" (Y) columns.
# ------------------------------------------------------------------
$wsFC = $wb.Worksheets.Item("FC_WEGE3")
$wsFC.Activate()

$wsFC.Columns("V:V").Insert(-4161, 0)

$wsFC.Range("V1").Value = "CAPEX"

$wsFC.Range("V2").Formula = "=X2+(Y2)"
$wsFC.Range("V3:V10").Formula = "=X3+(Y3)"

# Match the column width used by the neighbouring column (U).
$wsFC.Columns("V:V").ColumnWidth = 19.83

$wsFC.Range("Y9:Z9").Select()

# ------------------------------------------------------------------
# CGO_WEGE3: selection moves from A2:B2 to B5:B9 (and it is no
# longer the tab that is active/selected).
# ------------------------------------------------------------------
$wsCGO = $wb.Worksheets.Item("CGO_WEGE3")
$wsCGO.Activate()
$wsCGO.Range("B5:B9").Select()

# ------------------------------------------------------------------
# DRE_WEGE3: selection moves from A1:B10 to A1:XFD6 (whole rows
# 1-6 selected), and this becomes the active / selected sheet.
# ------------------------------------------------------------------
$wsDRE = $wb.Worksheets.Item("DRE_WEGE3")
$wsDRE.Activate()
$wsDRE.Range("A1:XFD6").Select()

Write-Host "edit complete"
